# Insert a new price-record row at row 83 ("Poroto verde" / Feria Lagunitas
# de Puerto Montt), pushing the existing rows 83-157 down to 84-158.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 83 — this shifts rows
# 83:157 down to 84:158 and grows the sheet's used range to A1:R158.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record's data.
$ws.Cells.Item(83, 1).Value = 4
$ws.Cells.Item(83, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(83, 3).Value = "Los Lagos"
$ws.Cells.Item(83, 4).Value = 45118
$ws.Cells.Item(83, 5).Value = 10
$ws.Cells.Item(83, 6).Value = 100112031
$ws.Cells.Item(83, 7).Value = "Poroto verde"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 45
$ws.Cells.Item(83, 11).Value = 31000
$ws.Cells.Item(83, 12).Value = 31000
$ws.Cells.Item(83, 13).Value = 31000
$ws.Cells.Item(83, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(83, 15).Value = "Perú"
$ws.Cells.Item(83, 16).Value = 1240
$ws.Cells.Item(83, 17).Value = 25
$ws.Cells.Item(83, 18).Value = "Hortaliza"
